$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "Pre Experimental Phase" (column C) self-report outcome values
$ws.Range("C2").Value = "Moderately stressful"
$ws.Range("C3").Value = "A little stressful"
$ws.Range("C4").Value = "A little stressful"
$ws.Range("C5").Value = "A little stressful"
$ws.Range("C6").Value = "Very stressful"
$ws.Range("C7").Value = "Very stressful"

# Update the selected cell to C8
$ws.Range("C8").Select()
